# Split the run "{m" into two runs "{" and "m", and split the run
# "someCustomService()}" into two runs "someCustomService()" and "}".
#
# The Word object model has no direct "split a run" verb, but adding a
# (temporary) bookmark at a position inside a run forces the run to be
# broken in two at that position; deleting the bookmark right afterwards
# removes the bookmark markers while leaving the now-separate runs in
# place (with no extra formatting residue).

$d = $word.ActiveDocument

# Locate the paragraph that holds the "{m:self.someCustomService()}" field.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("{m:self.someCustomService()}")) {
        $target = $p
    }
}

$start = $target.Range.Start
$end = $target.Range.End

# Split "{m" -> "{" | "m" : break right after the "{" (1 char in).
$splitAfterBrace = $d.Range($start + 1, $start + 1)
$d.Bookmarks.Add("m2docSplitA", $splitAfterBrace)
$d.Bookmarks.Item("m2docSplitA").Delete()

# Split "someCustomService()}" -> "someCustomService()" | "}" : break right
# before the closing "}", which is the next-to-last character of the
# paragraph (the very last position is the paragraph mark itself).
$splitBeforeBrace = $d.Range($end - 2, $end - 2)
$d.Bookmarks.Add("m2docSplitB", $splitBeforeBrace)
$d.Bookmarks.Item("m2docSplitB").Delete()
